$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-7 (Adam28-Itga4 LR pair, Neutrophils sending cluster)

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4793446666666667
$ws.Range("H2").Value = 1.438034
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.105124
$ws.Range("N2").Value = 3.315372
$ws.Range("O2").Value = 0.006910839970832482
$ws.Range("P2").Value = 0.006910839970832482
$ws.Range("Q2").Value = 0.5297352954053334
$ws.Range("R2").Value = 4.767617658648
$ws.Range("S2").Value = 0.006910839970832482
$ws.Range("T2").Value = 0.006910839970832482

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4793446666666667
$ws.Range("H3").Value = 1.438034
$ws.Range("O3").Value = 0.0002777950170396876
$ws.Range("P3").Value = 0.0002777950170396876
$ws.Range("Q3").Value = 0.02129376834577778
$ws.Range("R3").Value = 0.191643915112
$ws.Range("S3").Value = 0.0002777950170396876
$ws.Range("T3").Value = 0.0002777950170396876

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4793446666666667
$ws.Range("H4").Value = 1.438034
$ws.Range("M4").Value = 56.54517366666666
$ws.Range("N4").Value = 169.635521
$ws.Range("O4").Value = 0.3536025335919447
$ws.Range("P4").Value = 0.3536025335919447
$ws.Range("Q4").Value = 27.10462742285711
$ws.Range("R4").Value = 243.941646805714
$ws.Range("S4").Value = 0.3536025335919447
$ws.Range("T4").Value = 0.3536025335919447

# Row 5 (Target cluster: MuSCs)
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4793446666666667
$ws.Range("H5").Value = 1.438034
$ws.Range("M5").Value = 0.8044289999999998
$ws.Range("N5").Value = 2.413287
$ws.Range("O5").Value = 0.005030458199167516
$ws.Range("P5").Value = 0.005030458199167516
$ws.Range("Q5").Value = 0.385598750862
$ws.Range("R5").Value = 3.470388757757999
$ws.Range("S5").Value = 0.005030458199167516
$ws.Range("T5").Value = 0.005030458199167516

# Row 6 (Target cluster: Neutrophils)
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4793446666666667
$ws.Range("H6").Value = 1.438034
$ws.Range("M6").Value = 78.08909333333334
$ws.Range("N6").Value = 234.26728
$ws.Range("O6").Value = 0.4883264027331488
$ws.Range("P6").Value = 0.4883264027331488
$ws.Range("Q6").Value = 37.43159041416889
$ws.Range("R6").Value = 336.88431372752
$ws.Range("S6").Value = 0.4883264027331488
$ws.Range("T6").Value = 0.4883264027331488

# Row 7 (Target cluster: Resolving-Mac)
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4793446666666667
$ws.Range("H7").Value = 1.438034
$ws.Range("M7").Value = 23.323433
$ws.Range("N7").Value = 69.970299
$ws.Range("O7").Value = 0.1458519704878668
$ws.Range("P7").Value = 0.1458519704878668
$ws.Range("Q7").Value = 11.17996321690733
$ws.Range("R7").Value = 100.619668952166
$ws.Range("S7").Value = 0.1458519704878668
$ws.Range("T7").Value = 0.1458519704878668
